# Generate Report for Handoff:
#  - Status changes from "Handed back: in sync with en-US" to "Ready for handoff"
#    on the Overview sheet (E2/F2) and on each language sheet's Status cell (C2).
#  - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
#    are refreshed to reflect the new handoff.
#  - The (now shorter) "Status" columns are narrowed to fit the new text.
#    NB: ColumnWidth of 16.3826548258464 is the COM input that this engine's
#    pixel-snapping resolves to the narrow width used by the report template.
$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-15 10:56:55"
$wsOverview.Columns.Item(5).ColumnWidth = 16.3826548258464
$wsOverview.Columns.Item(6).ColumnWidth = 16.3826548258464

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-15 10:56:50"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3826548258464

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-15 10:56:55"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3826548258464
